$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7 (state "LITERAL") and a new column before
# column I (trigger "'" -> LITERAL). Excel auto-shifts existing cells and
# adjusts formulas, mirroring the diff's B7->B8, B8->B9, I->J shift.
$ws.Rows("7:7").Insert()
$ws.Columns("I:I").Insert()

# New column header I2: the literal-quote trigger character.
# B7 must be populated before I2 so the shared-string table gets
# "LITERAL" at index 13 and "'" at index 14, matching the target file.
$ws.Range("B7").Value = "LITERAL"

# Row 7 transitions: from state LITERAL, almost everything loops back to
# LITERAL except the trigger char in I (the closing quote) -> INIT.
$ws.Range("C7").Formula = "=B7"
$ws.Range("D7").Formula = "=B7"
$ws.Range("E7").Formula = "=B7"
$ws.Range("F7").Formula = "=B7"
$ws.Range("G7").Formula = "=B7"
$ws.Range("H7").Formula = "=B7"
$ws.Range("I7").Formula = "=B3"
$ws.Range("J7").Formula = "=B7"

# New column I: the trigger is a single quote character. Entering "''"
# makes Excel consume the first quote as the quote-prefix marker and
# store the literal single quote as the cell text (same pattern already
# used for the ':' / '-' / '_' header cells).
$ws.Range("I2").Value = "''"

# From SPECIAL_CHAR (row6) and LITERAL (row9 after shift) states, hitting
# a quote transitions into the LITERAL state (B7).
$ws.Range("I6").Formula = "=B7"
$ws.Range("I9").Formula = "=B7"

# Selection moves to I4 per the diff's sheetView.
$null = $ws.Range("I4").Select()
